# fix all function, add import cylinder
$wb = $excel.ActiveWorkbook

# --- Sheet "Create Cylinder" (sheet2): update test codes 19/08 -> 22/08 ---
$wsCyl = $wb.Worksheets.Item("Create Cylinder")
$wsCyl.Activate()
$wsCyl.Range("A2").Value = "TEST2208001"
$wsCyl.Range("A3").Value = "TEST2208002"
$wsCyl.Range("A4").Value = "TEST2208003"
$wsCyl.Range("A2:A4").Select()

# --- Sheet "Create New Account" (sheet1): replace account #5 with account #6,
#     and add a new account #8 row with a hyperlink ---
$wsAcc = $wb.Worksheets.Item("Create New Account")
$wsAcc.Activate()

# Row 2: account #6 (overwrite previous account #5 data)
$wsAcc.Range("A2").Value = "dailyphanphoi180806@gmail.com"
$wsAcc.Range("B2").Value = "Đại lý phân phối 6"
$wsAcc.Range("C2").Value = "DLPP180806"
$wsAcc.Range("F2").Value = "CDLPP180806"
$wsAcc.Range("G2").Value = "Chi nhánh DLPP 18/08 06"
$wsAcc.Range("H2").Value = "Đồng Nai"

# Re-point the existing hyperlink on A2 to the new account email
$wsAcc.Range("A2").Hyperlinks.Delete()
$wsAcc.Hyperlinks.Add($wsAcc.Range("A2"), "mailto:dailyphanphoi180806@gmail.com") | Out-Null

# Row 3: new account #8 - copy row 2's formatting first, then fill in the values
$wsAcc.Range("B2:H2").Copy() | Out-Null
$wsAcc.Range("B3:H3").PasteSpecial(-4122) | Out-Null
$wsAcc.Range("A2").Copy() | Out-Null
$wsAcc.Range("A3").PasteSpecial(-4122) | Out-Null

$wsAcc.Range("G3").Value = "Chi nhánh DLPP 18/08 08"
$wsAcc.Range("F3").Value = "CDLPP180808"
$wsAcc.Range("C3").Value = "DLPP180808"
$wsAcc.Range("D3").Value = "Quận 9"
$wsAcc.Range("E3").Value = "Nhóm thực tế"
$wsAcc.Range("B3").Value = "Đại lý phân phối 8"
$wsAcc.Range("H3").Value = "Đồng Nai"
$wsAcc.Range("A3").Value = "dailyphanphoi180808@gmail.com"

# Hyperlink for the new account email (A3), matching the existing A2 hyperlink pattern
$wsAcc.Hyperlinks.Add($wsAcc.Range("A3"), "mailto:dailyphanphoi180808@gmail.com") | Out-Null

$wsAcc.Range("A3").Select()
